$d = $word.ActiveDocument

# --- Title: "Default Title" -> "Sinh học" ---
# Single-run paragraph, safe to use Find/Replace.
$d.Content.Find.Execute("Default Title", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sinh học", 2)

# --- Question paragraphs: prefix "Câu N. " and consume the existing leading
# space (so we don't end up with a double space nor touch sibling runs'
# xml:space handling). We insert right at each paragraph's start boundary:
# the paragraph text already begins with a single space before the first
# word, so inserting "Câu N." (no trailing space) right there reuses that
# existing space as the separator and only rewrites the first <w:t> of the
# run, leaving sibling A/B/C/D runs completely untouched. ---

$captions = @(
    "Câu 1.",
    "Câu 2.",
    "Câu 3.",
    "Câu 4.",
    "Câu 5.",
    "Câu 6.",
    "Câu 7.",
    "Câu 8.",
    "Câu 9."
)

# Recompute each paragraph's start just before inserting, since earlier
# insertions shift later character offsets.
for ($q = 0; $q -lt $captions.Length; $q++) {
    $paraIndex = 2 + $q   # question paragraphs are Word paragraphs 2..10
    $start = $d.Paragraphs($paraIndex).Range.Start
    $ins = $d.Range($start, $start)
    $ins.InsertAfter($captions[$q])
}

# --- Answer key: drop the recorded letter answers for items 1-4 ---
$d.Content.Find.Execute("1. A ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1. ", 2)
$d.Content.Find.Execute("2. B ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2. ", 2)
$d.Content.Find.Execute("3. C ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3. ", 2)
$d.Content.Find.Execute("4. D ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4. ", 2)
